$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rotation ("E" column) updates: stored as plain numbers ---
$ws.Range("E37").Value = 180
$ws.Range("E49").Value = 90
$ws.Range("E53").Value = 270
$ws.Range("E54").Value = 0
$ws.Range("E55").Value = 90
$ws.Range("E56").Value = "90.000000"
$ws.Range("E57").Value = 90
$ws.Range("E58").Value = 180
$ws.Range("E59").Value = 180
$ws.Range("E60").Value = 90
$ws.Range("E61").Value = 90

# --- SW1 Mid Y update ---
$ws.Range("C101").Value = -13.4

# --- SW2 / SW3 Mid X updates ---
$ws.Range("B102").Value = 90.5
$ws.Range("B103").Value = 90.5

# --- U1..U11 Rotation updates ---
$ws.Range("E104").Value = "270.0"
$ws.Range("E105").Value = 90
$ws.Range("E106").Value = 270
$ws.Range("E107").Value = 180
$ws.Range("E108").Value = 180
$ws.Range("E109").Value = 180
$ws.Range("E110").Value = 0
$ws.Range("E111").Value = 270
$ws.Range("E112").Value = 90
$ws.Range("E113").Value = 180
$ws.Range("E114").Value = 0

# --- sheet view: scrolled + selection changed ---
$ws.Application.ActiveWindow.ScrollRow = 79
$ws.Range("B103").Select()
